$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.37041050194737
$ws.Range("D2").Value = 4.155204696756664
$ws.Range("E2").Value = 13.72610928902244
$ws.Range("F2").Value = 26.5227943986173
$ws.Range("G2").Value = 34.09956198786416
$ws.Range("H2").Value = 14.5584433102773
$ws.Range("L2").Value = 9.145448808988293
$ws.Range("O2").Value = 23.13022551420615
$ws.Range("C3").Value = 13.2699848469067
$ws.Range("D3").Value = 4.164241322338191
$ws.Range("E3").Value = 13.66501597118275
$ws.Range("F3").Value = 26.12679731536191
$ws.Range("G3").Value = 33.26842763544546
$ws.Range("H3").Value = 14.50374877142434
$ws.Range("L3").Value = 9.149145858536437
$ws.Range("O3").Value = 22.89358244757697
$ws.Range("C4").Value = 13.21156305405435
$ws.Range("D4").Value = 4.169974493241562
$ws.Range("E4").Value = 13.63080591497898
$ws.Range("F4").Value = 25.88873046564299
$ws.Range("G4").Value = 32.75873017636831
$ws.Range("H4").Value = 14.47347989785779
$ws.Range("L4").Value = 9.153051092502331
$ws.Range("O4").Value = 22.75360380579048
$ws.Range("C5").Value = 13.18859190932554
$ws.Range("D5").Value = 4.172357565646865
$ws.Range("E5").Value = 13.61770494133451
$ws.Range("F5").Value = 25.79312095759185
$ws.Range("G5").Value = 32.55152567640879
$ws.Range("H5").Value = 14.46198665716695
$ws.Range("L5").Value = 9.155053650312563
$ws.Range("O5").Value = 22.69796028144013
$ws.Range("C6").Value = 13.18482866954961
$ws.Range("D6").Value = 4.172756106706081
$ws.Range("E6").Value = 13.61558052323217
$ws.Range("F6").Value = 25.77733359755572
$ws.Range("G6").Value = 32.51716006840246
$ws.Range("H6").Value = 14.46012926374597
$ws.Range("L6").Value = 9.155411001309934
$ws.Range("O6").Value = 22.68880692487073
$ws.Range("C7").Value = 13.21124984418015
$ws.Range("D7").Value = 4.170006442449698
$ws.Range("E7").Value = 13.63062581798028
$ws.Range("F7").Value = 25.88743518790134
$ws.Range("G7").Value = 32.75593325152887
$ws.Range("H7").Value = 14.47332147840508
$ws.Range("L7").Value = 9.153076435169229
$ws.Range("O7").Value = 22.75284763561879
$ws.Range("C8").Value = 13.33512466194463
$ws.Range("D8").Value = 4.158282400405041
$ws.Range("E8").Value = 13.70436527329523
$ws.Range("F8").Value = 26.38527843459258
$ws.Range("G8").Value = 33.81306139919641
$ws.Range("H8").Value = 14.53890214088211
$ws.Range("L8").Value = 9.146384169886678
$ws.Range("O8").Value = 23.04756400656907
$ws.Range("C9").Value = 13.60273188941013
$ws.Range("D9").Value = 4.136741472071797
$ws.Range("E9").Value = 13.87468339608514
$ws.Range("F9").Value = 27.39610705352773
$ws.Range("G9").Value = 35.87645284445384
$ws.Range("H9").Value = 14.69338825875395
$ws.Range("L9").Value = 9.146232988110754
$ws.Range("O9").Value = 23.66491092064372
$ws.Range("C10").Value = 13.81293516011305
$ws.Range("D10").Value = 4.121777607462654
$ws.Range("E10").Value = 14.01477531325477
$ws.Range("F10").Value = 28.15212020161776
$ws.Range("G10").Value = 37.36806833338208
$ws.Range("H10").Value = 14.82205274249361
$ws.Range("L10").Value = 9.15401993156674
$ws.Range("O10").Value = 24.13857679840657
$ws.Range("C11").Value = 13.91117940161945
$ws.Range("D11").Value = 4.115152714106184
$ws.Range("E11").Value = 14.08158059664056
$ws.Range("F11").Value = 28.49733257179048
$ws.Range("G11").Value = 38.03778817837411
$ws.Range("H11").Value = 14.88372749745573
$ws.Range("L11").Value = 9.159272117890168
$ws.Range("O11").Value = 24.35755877411361
$ws.Range("C12").Value = 13.94872960505776
$ws.Range("D12").Value = 4.112669886466444
$ws.Range("E12").Value = 14.10730431498534
$ws.Range("F12").Value = 28.62811052037356
$ws.Range("G12").Value = 38.28985415784098
$ws.Range("H12").Value = 14.90752024431426
$ws.Range("L12").Value = 9.161506106580983
$ws.Range("O12").Value = 24.44091257073005
$ws.Range("C13").Value = 13.94062753326044
$ws.Range("D13").Value = 4.113203462258356
$ws.Range("E13").Value = 14.10174557542715
$ws.Range("F13").Value = 28.59994484855992
$ws.Range("G13").Value = 38.23563997656873
$ws.Range("H13").Value = 14.90237681339962
$ws.Range("L13").Value = 9.16101409038983
$ws.Range("O13").Value = 24.42294289741475
$ws.Range("C14").Value = 13.91426189276561
$ws.Range("D14").Value = 4.114947933691749
$ws.Range("E14").Value = 14.08368846160816
$ws.Range("F14").Value = 28.50809142899396
$ws.Range("G14").Value = 38.05855820458878
$ws.Range("H14").Value = 14.88567625069121
$ws.Range("L14").Value = 9.159451002355386
$ws.Range("O14").Value = 24.3644081484729
$ws.Range("C15").Value = 13.89815650436399
$ws.Range("D15").Value = 4.116019832769152
$ws.Range("E15").Value = 14.07268293517138
$ws.Range("F15").Value = 28.45183163900171
$ws.Range("G15").Value = 37.94988186959795
$ws.Range("H15").Value = 14.87550325658975
$ws.Range("L15").Value = 9.158525459155161
$ws.Range("O15").Value = 24.32860772178099
$ws.Range("C16").Value = 13.80656489398881
$ws.Range("D16").Value = 4.12221419876768
$ws.Range("E16").Value = 14.0104699580445
$ws.Range("F16").Value = 28.12957441757662
$ws.Range("G16").Value = 37.32409916418757
$ws.Range("H16").Value = 14.81808422271523
$ws.Range("L16").Value = 9.153711033056513
$ws.Range("O16").Value = 24.12433005285329
$ws.Range("C17").Value = 13.751026925196
$ws.Range("D17").Value = 4.126060683422858
$ws.Range("E17").Value = 13.97308092542883
$ws.Range("F17").Value = 27.93211667036691
$ws.Range("G17").Value = 36.93774073662581
$ws.Range("H17").Value = 14.78365497450513
$ws.Range("L17").Value = 9.151194975731721
$ws.Range("O17").Value = 23.99985807693579
$ws.Range("C18").Value = 13.71933155048398
$ws.Range("D18").Value = 4.128290253099933
$ws.Range("E18").Value = 13.95186638201745
$ws.Range("F18").Value = 27.81867415254598
$ws.Range("G18").Value = 36.71470216004705
$ws.Range("H18").Value = 14.76414908063748
$ws.Range("L18").Value = 9.149908803807708
$ws.Range("O18").Value = 23.92859985185313
$ws.Range("C19").Value = 13.70864363409635
$ws.Range("D19").Value = 4.129048106701058
$ws.Range("E19").Value = 13.94473390647453
$ws.Range("F19").Value = 27.78029090189899
$ws.Range("G19").Value = 36.63905374850774
$ws.Range("H19").Value = 14.7575961366013
$ws.Range("L19").Value = 9.149501000575132
$ws.Range("O19").Value = 23.90453280603566
$ws.Range("C20").Value = 13.75691352165736
$ws.Range("D20").Value = 4.125649443604354
$ws.Range("E20").Value = 13.97703108809826
$ws.Range("F20").Value = 27.95312391127026
$ws.Range("G20").Value = 36.97895569104341
$ws.Range("H20").Value = 14.78728940031411
$ws.Range("L20").Value = 9.1514461568075
$ws.Range("O20").Value = 24.01307422436246
$ws.Range("C21").Value = 13.92199693653995
$ws.Range("D21").Value = 4.11443484052238
$ws.Range("E21").Value = 14.08898084690301
$ws.Range("F21").Value = 28.5350705948232
$ws.Range("G21").Value = 38.11061535280589
$ws.Range("H21").Value = 14.8905698456515
$ws.Range("L21").Value = 9.159903474271859
$ws.Range("O21").Value = 24.38159014469169
$ws.Range("C22").Value = 14.03189830191671
$ws.Range("D22").Value = 4.107256118724015
$ws.Range("E22").Value = 14.1646214060468
$ws.Range("F22").Value = 28.91565940356292
$ws.Range("G22").Value = 38.84111764272442
$ws.Range("H22").Value = 14.96061455826215
$ws.Range("L22").Value = 9.166858972503363
$ws.Range("O22").Value = 24.62490897187235
$ws.Range("C23").Value = 13.97306795110978
$ws.Range("D23").Value = 4.11107386073448
$ws.Range("E23").Value = 14.12402979452889
$ws.Range("F23").Value = 28.71255181703465
$ws.Range("G23").Value = 38.45215282951923
$ws.Range("H23").Value = 14.92300244718494
$ws.Range("L23").Value = 9.163016315270008
$ws.Range("O23").Value = 24.49484349743259
$ws.Range("C24").Value = 13.75425145984127
$ws.Range("D24").Value = 4.125835308475732
$ws.Range("E24").Value = 13.9752443433438
$ws.Range("F24").Value = 27.94362628268639
$ws.Range("G24").Value = 36.96032524213592
$ws.Range("H24").Value = 14.78564537839155
$ws.Range("L24").Value = 9.15133209834373
$ws.Range("O24").Value = 24.0070982565469
$ws.Range("C25").Value = 13.52783970910518
$ws.Range("D25").Value = 4.142415879446895
$ws.Range("E25").Value = 13.82592466555292
$ws.Range("F25").Value = 27.11974928405034
$ws.Range("G25").Value = 35.3212553033973
$ws.Range("H25").Value = 14.64888441999574
$ws.Range("L25").Value = 9.14488554871282
$ws.Range("O25").Value = 23.49408141794444
